$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo in dish name
$ws.Range("D9").Value = "Чизу Рамен"

# Add counter values (orm query counters fix)
$ws.Range("G7").Value = 59
$ws.Range("G8").Value = 31

# Move selection/cursor
$ws.Range("I6").Select()
